$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32: Automata for the People | Crab Oil
$ws.Range("H32").Value = 2507.8572
$ws.Range("I32").Value = 2018.3334
$ws.Range("K32").Value = 2018.3334
$ws.Range("M32").Value = -1692.3334
# Row 53: No Accounting for Waste | Enchanted Electrum Ink
$ws.Range("H53").Value = 423.13794
$ws.Range("I53").Value = 305.57144
$ws.Range("J53").Value = 532.86664
$ws.Range("K53").Value = 305.57144
$ws.Range("L53").Value = 532.86664
$ws.Range("M53").Value = 331.42856
$ws.Range("N53").Value = -1806.86664
# Row 70: Consecrating Congregation | Holy Water
$ws.Range("H70").Value = 5147.5
$ws.Range("I70").Value = 3500
$ws.Range("J70").Value = 5559.375
$ws.Range("K70").Value = 10500
$ws.Range("L70").Value = 16678.125
$ws.Range("M70").Value = -10230
$ws.Range("N70").Value = -17218.125
# Row 73: Curbing the Contagion (L) | Holy Water
$ws.Range("H73").Value = 5147.5
$ws.Range("I73").Value = 3500
$ws.Range("J73").Value = 5559.375
$ws.Range("K73").Value = 10500
$ws.Range("L73").Value = 16678.125
$ws.Range("M73").Value = -9564
$ws.Range("N73").Value = -18550.125
# Row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 3117.9524
$ws.Range("I98").Value = 1168.5385
$ws.Range("J98").Value = 6285.75
$ws.Range("K98").Value = 1168.5385
$ws.Range("L98").Value = 6285.75
$ws.Range("M98").Value = 329.4614999999999
$ws.Range("N98").Value = -9281.75
# Row 116: Growing Up | Growth Formula Kappa
$ws.Range("H116").Value = 308580.8
$ws.Range("I116").Value = 529029.8
$ws.Range("J116").Value = 9400
$ws.Range("K116").Value = 529029.8
$ws.Range("L116").Value = 9400
$ws.Range("M116").Value = -525587.8
$ws.Range("N116").Value = -16284
# Row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 3117.9524
$ws.Range("I122").Value = 1168.5385
$ws.Range("J122").Value = 6285.75
$ws.Range("K122").Value = 3505.6155
$ws.Range("L122").Value = 18857.25
$ws.Range("M122").Value = -1055.6155
$ws.Range("N122").Value = -23757.25
# Row 129: Practical Command | Commanding Craftsman's Draught
$ws.Range("H129").Value = 885.91
$ws.Range("J129").Value = 907.33685
$ws.Range("L129").Value = 2722.01055
$ws.Range("N129").Value = -12722.01055
# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2390.07
$ws.Range("I138").Value = 996.2
$ws.Range("J138").Value = 2636.047
$ws.Range("K138").Value = 2988.6
$ws.Range("L138").Value = 7908.141
$ws.Range("M138").Value = 2151.4
$ws.Range("N138").Value = -18188.141

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots | Bronze Ingot
$ws.Range("H2").Value = 549.975
$ws.Range("I2").Value = 493.27585
$ws.Range("K2").Value = 493.27585
$ws.Range("M2").Value = -380.27585
# Row 116: No Scope | Titanbronze Ingot
$ws.Range("H116").Value = 549.975
$ws.Range("I116").Value = 493.27585
$ws.Range("K116").Value = 493.27585
$ws.Range("M116").Value = 1800.72415
# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 2838.5454
$ws.Range("I122").Value = 1580.4445
$ws.Range("K122").Value = 4741.333500000001
$ws.Range("M122").Value = -2291.333500000001

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells | Bronze Ingot
$ws.Range("H3").Value = 549.975
$ws.Range("I3").Value = 493.27585
$ws.Range("K3").Value = 493.27585
$ws.Range("M3").Value = -379.27585
# Row 128: Mangalomania | Manganese Ingot
$ws.Range("H128").Value = 1225
$ws.Range("I128").Value = 1225
$ws.Range("K128").Value = 3675
$ws.Range("M128").Value = -1185

$ws = $wb.Worksheets.Item("CRP")
# Row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws.Range("H122").Value = 2871.0715
$ws.Range("I122").Value = 1932.9166
$ws.Range("J122").Value = 8500
$ws.Range("K122").Value = 5798.7498
$ws.Range("L122").Value = 25500
$ws.Range("M122").Value = -3348.7498
$ws.Range("N122").Value = -30400

$ws = $wb.Worksheets.Item("CUL")
# Row 17: Chew the Fat | Grilled Dodo
$ws.Range("H17").Value = 4190
$ws.Range("I17").Value = 400
$ws.Range("J17").Value = 7980
$ws.Range("K17").Value = 1200
$ws.Range("L17").Value = 23940
$ws.Range("M17").Value = -1031
$ws.Range("N17").Value = -24278
# Row 23: Sweet Smell of Success | Lavender Oil
$ws.Range("H23").Value = 238
$ws.Range("I23").Value = 50
$ws.Range("J23").Value = 256.8
$ws.Range("K23").Value = 150
$ws.Range("L23").Value = 770.4000000000001
$ws.Range("M23").Value = 85
$ws.Range("N23").Value = -1240.4
# Row 40: True Grits | Cornmeal
$ws.Range("H40").Value = 500
$ws.Range("I40").Value = 201
$ws.Range("J40").Value = 799
$ws.Range("K40").Value = 804
$ws.Range("L40").Value = 3196
$ws.Range("M40").Value = -735
$ws.Range("N40").Value = -3334
# Row 44: No More Dumpster Diving | Knight's Bread
$ws.Range("H44").Value = 714.25
$ws.Range("I44").Value = 329
$ws.Range("J44").Value = 1099.5
$ws.Range("K44").Value = 987
$ws.Range("L44").Value = 3298.5
$ws.Range("M44").Value = -589
$ws.Range("N44").Value = -4094.5
# Row 70: Persona non Gratin | Dhalmel Gratin
$ws.Range("H70").Value = 2585.1667
$ws.Range("J70").Value = 2999.75
$ws.Range("L70").Value = 8999.25
$ws.Range("N70").Value = -9629.25
# Row 73: Recipe for Disaster (L) | Dhalmel Gratin
$ws.Range("H73").Value = 2585.1667
$ws.Range("J73").Value = 2999.75
$ws.Range("L73").Value = 8999.25
$ws.Range("N73").Value = -11183.25
# Row 121: A Cookie for Your Troubles | Coffee Biscuit
$ws.Range("H121").Value = 2681.9106
$ws.Range("I121").Value = 200
$ws.Range("J121").Value = 2727.0364
$ws.Range("K121").Value = 600
$ws.Range("L121").Value = 8181.1092
$ws.Range("M121").Value = 710
$ws.Range("N121").Value = -10801.1092
# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 7463596.5
$ws.Range("I131").Value = 125000424
$ws.Range("J131").Value = 940.3016
$ws.Range("K131").Value = 375001272
$ws.Range("L131").Value = 2820.9048
$ws.Range("M131").Value = -374996232
$ws.Range("N131").Value = -12900.9048

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value = 22729956
$ws.Range("I80").Value = 41668816
$ws.Range("K80").Value = 41668816
$ws.Range("M80").Value = -41667818
# Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value = 22729956
$ws.Range("I83").Value = 41668816
$ws.Range("K83").Value = 208344080
$ws.Range("M83").Value = -208339088

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 3790.5833
$ws.Range("I7").Value = 2476.3333
$ws.Range("J7").Value = 7733.3335
$ws.Range("K7").Value = 2476.3333
$ws.Range("L7").Value = 7733.3335
$ws.Range("M7").Value = -2364.3333
$ws.Range("N7").Value = -7957.3335
# Row 16: Saddle Sore | Hard Leather
$ws.Range("H16").Value = 1177.8182
$ws.Range("I16").Value = 995.6
$ws.Range("K16").Value = 995.6
$ws.Range("M16").Value = -825.6
# Row 40: Best Served Toad | Toad Leather
$ws.Range("H40").Value = 8843.65
$ws.Range("I40").Value = 11160.5
$ws.Range("J40").Value = 7850.7144
$ws.Range("K40").Value = 11160.5
$ws.Range("L40").Value = 7850.7144
$ws.Range("M40").Value = -11024.5
$ws.Range("N40").Value = -8122.7144
# Row 46: Supply Side Logic | Boar Leather
$ws.Range("H46").Value = 3608.5
$ws.Range("J46").Value = 3617.3333
$ws.Range("L46").Value = 3617.3333
$ws.Range("N46").Value = -3993.3333
# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 4291.923
$ws.Range("I122").Value = 2179.5
$ws.Range("J122").Value = 11333.333
$ws.Range("K122").Value = 6538.5
$ws.Range("L122").Value = 33999.999
$ws.Range("M122").Value = -4088.5
$ws.Range("N122").Value = -38899.999
# Row 123: Running up the Tabi | Gajaskin Tabi
$ws.Range("H123").Value = 50429
$ws.Range("J123").Value = 50429
$ws.Range("L123").Value = 50429
$ws.Range("N123").Value = -60229
# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 3790.5833
$ws.Range("I126").Value = 2476.3333
$ws.Range("J126").Value = 7733.3335
$ws.Range("K126").Value = 7428.999899999999
$ws.Range("L126").Value = 23200.0005
$ws.Range("M126").Value = -4958.999899999999
$ws.Range("N126").Value = -28140.0005
